$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$matchName = "CDF T7 VS Pierrelatte (R3)"
$period = "Global"

$startRow = 858
$endRow = 870

# Row data: E (player), F (poste), G (temps joue), H..V (numeric stats)
$rows = @(
    @{ E="Naim Dhib";       F="center midfield"; G="01:40:47"; H=10;    I=1.61; J=8.3699999999999992;  K=1.1399999999999999;  L=0.36;                 M=0.13; N=0;    O=8;  P=5.95;                Q=29.73; R=4.5;                  S=47; T=5;  U=39; V=6  }
    @{ E="Maé Clavel";      F="left back";       G="00:48:08"; H=4.96;  I=0.73; J=4.22;               K=0.49;                 L=0.23;                 M=0.01; N=0;    O=2;  P=6.18;                Q=25.83; R=4.47;                 S=22; T=5;  U=15; V=5  }
    @{ E="Karim Belmahi";   F="left forward";    G="01:14:20"; H=8.32;  I=1.76; J=6.53;               K=1.0900000000000001;  L=0.54;                 M=0.16; N=0;    O=14; P=6.62;                Q=29.84; R=4.6399999999999997;   S=50; T=11; U=37; V=14 }
    @{ E="Levy Ndoutoume";  F="left back";       G="01:00:45"; H=6.67;  I=1.25; J=5.4;                K=0.82;                 L=0.38;                 M=0.06; N=0;    O=7;  P=6.58;                Q=27.36; R=4.82;                 S=33; T=3;  U=22; V=8  }
    @{ E="Ilyes Boughanmi"; F="center forward";  G="01:14:04"; H=6.65;  I=0.84; J=5.8;                K=0.62;                 L=0.2;                  M=0.03; N=0;    O=3;  P=5.35;                Q=27.68; R=5.0199999999999996;   S=29; T=3;  U=25; V=8  }
    @{ E="Sofiane Belle";   F="left forward";    G="00:24:57"; H=2.4;   I=0.4;  J=1.99;               K=0.31;                 L=0.06;                 M=0.03; N=0;    O=2;  P=5.82;                Q=28.88; R=4.26;                 S=6;  T=1;  U=6;  V=1  }
    @{ E="Yoan Zouma";      F="center back";     G="01:41:04"; H=8.59;  I=0.85; J=7.73;               K=0.71;                 L=0.14000000000000001; M=0.01; N=0;    O=1;  P=5.1100000000000003;  Q=26.24; R=4.42;                 S=21; T=1;  U=18; V=0  }
    @{ E="Yoann Martelat";  F="center midfield"; G="01:40:48"; H=12.24; I=2.2599999999999998; J=9.9600000000000009; K=1.81;    L=0.44;                 M=0.03; N=0;    O=3;  P=7.24;                Q=26.47; R=4.37;                 S=30; T=2;  U=17; V=3  }
    @{ E="Naim Ighbane";    F="center back";     G="01:14:28"; H=7.11;  I=0.89; J=6.21;               K=0.57999999999999996; L=0.24;                 M=0.08; N=0;    O=6;  P=5.51;                Q=27.36; R=4.1100000000000003;   S=19; T=3;  U=10; V=4  }
    @{ E="Mattheo Haon";    F="right back";      G="01:40:56"; H=11.11; I=1.99; J=9.11;               K=1.33;                 L=0.5;                  M=0.18; N=0;    O=14; P=6.56;                Q=30.03; R=5.09;                 S=40; T=9;  U=38; V=16 }
    @{ E="Emmanuel Valey";  F="left forward";    G="00:13:26"; H=1.73;  I=0.47; J=1.25;               K=0.25;                 L=0.16;                 M=0.07; N=0;    O=6;  P=7.7;                 Q=28.25; R=4.08;                 S=12; T=2;  U=9;  V=3  }
    @{ E="Amir Etien";      F="right forward";   G="01:26:15"; H=7.6;   I=1.39; J=6.19;               K=0.65;                 L=0.47;                 M=0.27; N=0.02; O=19; P=5.2;                 Q=31.98; R=5.65;                 S=42; T=11; U=23; V=15 }
    @{ E="Jeremie Laurent"; F="left forward";    G="00:15:30"; H=1.74;  I=0.49; J=1.24;               K=0.34;                 L=0.08;                 M=0.07; N=0;    O=5;  P=6.63;                Q=29.89; R=4.78;                 S=7;  T=5;  U=7;  V=4  }
)

# 1) Set column G (Temps joue) first for every new row, in row order, so the
#    new shared-string entries for the time values are created before the
#    match-name string (matches original authoring order).
$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 7).Value = $row.G
    $r++
}

# 2) Copy the date cell format from the last existing row down to the new rows
$ws.Cells.Item(857, 2).Copy()
$ws.Range("B$startRow`:B$endRow").PasteSpecial(-4122)

# 3) Fill in the remaining columns for each new row
$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $matchName
    $ws.Cells.Item($r, 2).Value = 45977
    $ws.Cells.Item($r, 3).Value = $period
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
    $r++
}

# Update the sheet scroll position / active selection to reflect the new data
$ws.Activate()
$ws.Range("E875").Select()
$excel.ActiveWindow.ScrollRow = 840
